$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.200.38'
$ws.Range('E2').Value = '  +2.08%  '

$ws.Range('D3').Value = '2.346.09'
$ws.Range('E3').Value = '  +6.04%  '

$ws.Range('E4').Value = '  -0.24%  '

$ws.Range('D5').Value = "'" + '313.27'
$ws.Range('E5').Value = '  +6.01%  '

$ws.Range('D6').Value = "'" + '108.98'
$ws.Range('E6').Value = '  +1.30%  '

$ws.Range('E7').Value = '  +3.38%  '

$ws.Range('E8').Value = '  -0.21%  '

$ws.Range('D9').Value = "'" + '0.636'
$ws.Range('E9').Value = '  +7.30%  '

$ws.Range('D10').Value = "'" + '42.87'
$ws.Range('E10').Value = '  -1.32%  '

$ws.Range('D11').Value = "'" + '0.0937'
$ws.Range('E11').Value = '  +3.35%  '

$ws.Range('D12').Value = "'" + '8.82'
$ws.Range('E12').Value = '  +1.02%  '

$ws.Range('E13').Value = '  +9.06%  '

$ws.Range('D14').Value = "'" + '0.105'
$ws.Range('E14').Value = '  +2.46%  '

$ws.Range('E15').Value = '  +9.25%  '

$ws.Range('D16').Value = '2.702.91'
$ws.Range('E16').Value = '  +6.08%  '

$ws.Range('D17').Value = '2.351.03'
$ws.Range('E17').Value = '  +5.48%  '

$ws.Range('D18').Value = '43.189.62'
$ws.Range('E18').Value = '  +2.18%  '

$ws.Range('D19').Value = "'" + '0.0000108'
$ws.Range('E19').Value = '  +3.42%  '

$ws.Range('D20').Value = "'" + '7.24'
$ws.Range('E20').Value = '  -1.78%  '

$ws.Range('D21').Value = "'" + '75.29'
$ws.Range('E21').Value = '  +3.87%  '

$ws.Range('D22').Value = "'" + '2.58'
$ws.Range('E22').Value = '  +12.33%  '

$ws.Range('E23').Value = '  -1.01%  '

$ws.Range('D24').Value = "'" + '254.69'
$ws.Range('E24').Value = '  +11.92%  '

$ws.Range('D25').Value = "'" + '9.08'
$ws.Range('E25').Value = '  +0.87%  '

$ws.Range('D26').Value = "'" + '12.00'
$ws.Range('E26').Value = '  +3.91%  '

$ws.Range('D27').Value = "'" + '0.999'
$ws.Range('E27').Value = '  +0.01%  '

$ws.Range('D28').Value = "'" + '39.08'
$ws.Range('E28').Value = '  +1.91%  '

$ws.Range('D29').Value = "'" + '2.25'
$ws.Range('E29').Value = '  +1.16%  '

$ws.Range('D30').Value = "'" + '22.31'
$ws.Range('E30').Value = '  +7.02%  '

$ws.Range('D31').Value = "'" + '173.76'
$ws.Range('E31').Value = '  +0.39%  '

$ws.Range('E32').Value = '  -0.79%  '

$ws.Range('D33').Value = "'" + '0.0925'
$ws.Range('E33').Value = '  +4.29%  '

$ws.Range('D34').Value = "'" + '6.07'

$ws.Range('D35').Value = "'" + '0.132'
$ws.Range('E35').Value = '  +5.75%  '

$ws.Range('D36').Value = "'" + '4.95'
$ws.Range('E36').Value = '  -2.11%  '

$ws.Range('D37').Value = "'" + '0.0377'
$ws.Range('E37').Value = '  +3.43%  '

$ws.Range('D38').Value = "'" + '4.12'
$ws.Range('E38').Value = '  -4.00%  '

$ws.Range('E39').Value = '  +2.55%  '

$ws.Range('D40').Value = "'" + '2.69'
$ws.Range('E40').Value = '  +11.30%  '

$ws.Range('D41').Value = "'" + '72.53'
$ws.Range('E41').Value = '  +1.67%  '

$ws.Range('E42').Value = '  +14.85%  '

$ws.Range('D43').Value = "'" + '0.233'
$ws.Range('E43').Value = '  +1.31%  '

$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = "'" + '1.00'
$ws.Range('E44').Value = '  +0.02%  '

$ws.Range('B45').Value = 'Celestia'
$ws.Range('C45').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D45').Value = "'" + '12.69'
$ws.Range('E45').Value = '  +1.60%  '

$ws.Range('D46').Value = "'" + '5.61'
$ws.Range('E46').Value = '  +3.82%  '

$ws.Range('D47').Value = "'" + '9.30'
$ws.Range('E47').Value = '  +11.91%  '

$ws.Range('D48').Value = "'" + '110.69'
$ws.Range('E48').Value = '  +7.38%  '

$ws.Range('E49').Value = '  -0.89%  '

$ws.Range('E50').Value = '  +4.16%  '

$ws.Range('D51').Value = "'" + '69.58'
$ws.Range('E51').Value = '  +5.14%  '
